$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Add the new "altimu_meas" sheet after the last existing sheet
# ---------------------------------------------------------------
$ws0 = $wb.Worksheets.Item(1)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "altimu_meas"

# ---------------------------------------------------------------
# 2) Header labels (row 2)
# ---------------------------------------------------------------
$ws.Range("A2").Value = "time_ms"
$ws.Range("B2").Value = "pressure"
$ws.Range("C2").Value = "temp"
$ws.Range("D2").Value = "vel_x"
$ws.Range("E2").Value = "accel_x"
$ws.Range("F2").Value = "vel_y"
$ws.Range("G2").Value = "accel_y"
$ws.Range("H2").Value = "vel_z"
$ws.Range("I2").Value = "accel_z"
$ws.Range("J2").Value = "mag_x"
$ws.Range("K2").Value = "mag_y"
$ws.Range("L2").Value = "mag_z"

# ---------------------------------------------------------------
# 3) Flight-data sample rows (3-24)
# ---------------------------------------------------------------
# row 3: A3=24712  B3=2213225  C3=3226  H3=3996  I3=7686  J3=64871  K3=906  L3=965
$ws.Range("A3").Value = 24712
$ws.Range("B3").Value = 2213225
$ws.Range("C3").Value = 3226
$ws.Range("H3").Value = 3996
$ws.Range("I3").Value = 7686
$ws.Range("J3").Value = 64871
$ws.Range("K3").Value = 906
$ws.Range("L3").Value = 965
# row 4: A4=25029  B4=2175746  C4=3211  D4=69  E4=63812  F4=582  G4=654  H4=4564  I4=7823  J4=64335  K4=502  L4=943
$ws.Range("A4").Value = 25029
$ws.Range("B4").Value = 2175746
$ws.Range("C4").Value = 3211
$ws.Range("D4").Value = 69
$ws.Range("E4").Value = 63812
$ws.Range("F4").Value = 582
$ws.Range("G4").Value = 654
$ws.Range("H4").Value = 4564
$ws.Range("I4").Value = 7823
$ws.Range("J4").Value = 64335
$ws.Range("K4").Value = 502
$ws.Range("L4").Value = 943
# row 5: A5=39917  D5=27  E5=65499  F5=65493  G5=21  H5=64884  I5=806
$ws.Range("A5").Value = 39917
$ws.Range("D5").Value = 27
$ws.Range("E5").Value = 65499
$ws.Range("F5").Value = 65493
$ws.Range("G5").Value = 21
$ws.Range("H5").Value = 64884
$ws.Range("I5").Value = 806
# row 6: A6=40237  B6=1115428  C6=3258  D6=65535  E6=65500  F6=44  G6=65423  H6=65018  I6=750  J6=63952  K6=924  L6=896
$ws.Range("A6").Value = 40237
$ws.Range("B6").Value = 1115428
$ws.Range("C6").Value = 3258
$ws.Range("D6").Value = 65535
$ws.Range("E6").Value = 65500
$ws.Range("F6").Value = 44
$ws.Range("G6").Value = 65423
$ws.Range("H6").Value = 65018
$ws.Range("I6").Value = 750
$ws.Range("J6").Value = 63952
$ws.Range("K6").Value = 924
$ws.Range("L6").Value = 896
# row 7: A7=55125  B7=819957  C7=3245
$ws.Range("A7").Value = 55125
$ws.Range("B7").Value = 819957
$ws.Range("C7").Value = 3245
# row 8: A8=55229  D8=65522  E8=65510  F8=65477  G8=65532  H8=1089  I8=164
$ws.Range("A8").Value = 55229
$ws.Range("D8").Value = 65522
$ws.Range("E8").Value = 65510
$ws.Range("F8").Value = 65477
$ws.Range("G8").Value = 65532
$ws.Range("H8").Value = 1089
$ws.Range("I8").Value = 164
# row 9: A9=55440  B9=818682  C9=3245  D9=65519  E9=65477  F9=65491  G9=65520  H9=1155  I9=159  J9=65174  K9=1232  L9=700
$ws.Range("A9").Value = 55440
$ws.Range("B9").Value = 818682
$ws.Range("C9").Value = 3245
$ws.Range("D9").Value = 65519
$ws.Range("E9").Value = 65477
$ws.Range("F9").Value = 65491
$ws.Range("G9").Value = 65520
$ws.Range("H9").Value = 1155
$ws.Range("I9").Value = 159
$ws.Range("J9").Value = 65174
$ws.Range("K9").Value = 1232
$ws.Range("L9").Value = 700
# row 10: A10=55546  D10=65516  E10=65506
$ws.Range("A10").Value = 55546
$ws.Range("D10").Value = 65516
$ws.Range("E10").Value = 65506
# row 11: A11=70327  B11=799316  C11=3141  D11=22  E11=15  F11=65475  G11=100  H11=649  I11=313  J11=64857  K11=1531  L11=65194
$ws.Range("A11").Value = 70327
$ws.Range("B11").Value = 799316
$ws.Range("C11").Value = 3141
$ws.Range("D11").Value = 22
$ws.Range("E11").Value = 15
$ws.Range("F11").Value = 65475
$ws.Range("G11").Value = 100
$ws.Range("H11").Value = 649
$ws.Range("I11").Value = 313
$ws.Range("J11").Value = 64857
$ws.Range("K11").Value = 1531
$ws.Range("L11").Value = 65194
# row 12: A12=85426  B12=1019289  C12=2806  F12=6  G12=65137  H12=1094  I12=905  J12=64876  K12=824  L12=64868
$ws.Range("A12").Value = 85426
$ws.Range("B12").Value = 1019289
$ws.Range("C12").Value = 2806
$ws.Range("F12").Value = 6
$ws.Range("G12").Value = 65137
$ws.Range("H12").Value = 1094
$ws.Range("I12").Value = 905
$ws.Range("J12").Value = 64876
$ws.Range("K12").Value = 824
$ws.Range("L12").Value = 64868
# row 13: A13=85532  B13=1016676  C13=2804  D13=65514  E13=65484  F13=2  G13=24  H13=1053  I13=940  J13=64862  K13=771  L13=64853
$ws.Range("A13").Value = 85532
$ws.Range("B13").Value = 1016676
$ws.Range("C13").Value = 2804
$ws.Range("D13").Value = 65514
$ws.Range("E13").Value = 65484
$ws.Range("F13").Value = 2
$ws.Range("G13").Value = 24
$ws.Range("H13").Value = 1053
$ws.Range("I13").Value = 940
$ws.Range("J13").Value = 64862
$ws.Range("K13").Value = 771
$ws.Range("L13").Value = 64853
# row 14: A14=85743  B14=1027103  C14=2803  D14=65508  E14=65447  F14=65526  G14=65272  H14=919  I14=905  J14=64823  K14=688  L14=64854
$ws.Range("A14").Value = 85743
$ws.Range("B14").Value = 1027103
$ws.Range("C14").Value = 2803
$ws.Range("D14").Value = 65508
$ws.Range("E14").Value = 65447
$ws.Range("F14").Value = 65526
$ws.Range("G14").Value = 65272
$ws.Range("H14").Value = 919
$ws.Range("I14").Value = 905
$ws.Range("J14").Value = 64823
$ws.Range("K14").Value = 688
$ws.Range("L14").Value = 64854
# row 15: A15=100738  B15=1487106  C15=2313  F15=17  G15=65202  H15=6  I15=1792  J15=64201  K15=1013  L15=64799
$ws.Range("A15").Value = 100738
$ws.Range("B15").Value = 1487106
$ws.Range("C15").Value = 2313
$ws.Range("F15").Value = 17
$ws.Range("G15").Value = 65202
$ws.Range("H15").Value = 6
$ws.Range("I15").Value = 1792
$ws.Range("J15").Value = 64201
$ws.Range("K15").Value = 1013
$ws.Range("L15").Value = 64799
# row 16: A16=100950  B16=1494342  C16=2308  D16=30  E16=79  F16=9  G16=65403  H16=116  I16=1892  J16=64220  K16=1003  L16=64664
$ws.Range("A16").Value = 100950
$ws.Range("B16").Value = 1494342
$ws.Range("C16").Value = 2308
$ws.Range("D16").Value = 30
$ws.Range("E16").Value = 79
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 65403
$ws.Range("H16").Value = 116
$ws.Range("I16").Value = 1892
$ws.Range("J16").Value = 64220
$ws.Range("K16").Value = 1003
$ws.Range("L16").Value = 64664
# row 17: A17=115945  B17=2161418  C17=2105  H17=64991  I17=2021  J17=64393  K17=1087  L17=64771
$ws.Range("A17").Value = 115945
$ws.Range("B17").Value = 2161418
$ws.Range("C17").Value = 2105
$ws.Range("H17").Value = 64991
$ws.Range("I17").Value = 2021
$ws.Range("J17").Value = 64393
$ws.Range("K17").Value = 1087
$ws.Range("L17").Value = 64771
# row 18: A18=116260  B18=2158255  C18=2106  D18=61  E18=64950  F18=65475  G18=65421  H18=64988  I18=2146  J18=64367  K18=1084  L18=64759
$ws.Range("A18").Value = 116260
$ws.Range("B18").Value = 2158255
$ws.Range("C18").Value = 2106
$ws.Range("D18").Value = 61
$ws.Range("E18").Value = 64950
$ws.Range("F18").Value = 65475
$ws.Range("G18").Value = 65421
$ws.Range("H18").Value = 64988
$ws.Range("I18").Value = 2146
$ws.Range("J18").Value = 64367
$ws.Range("K18").Value = 1084
$ws.Range("L18").Value = 64759
# row 19: A19=116366  D19=24  E19=64644  F19=65531  G19=31  H19=64817  I19=2256  J19=64359  K19=1070
$ws.Range("A19").Value = 116366
$ws.Range("D19").Value = 24
$ws.Range("E19").Value = 64644
$ws.Range("F19").Value = 65531
$ws.Range("G19").Value = 31
$ws.Range("H19").Value = 64817
$ws.Range("I19").Value = 2256
$ws.Range("J19").Value = 64359
$ws.Range("K19").Value = 1070
# row 20: A20=131253  B20=2977769  C20=2088  D20=57  E20=64640  F20=65535  G20=65244  H20=63759  I20=2266  J20=64641  K20=1118  L20=64793
$ws.Range("A20").Value = 131253
$ws.Range("B20").Value = 2977769
$ws.Range("C20").Value = 2088
$ws.Range("D20").Value = 57
$ws.Range("E20").Value = 64640
$ws.Range("F20").Value = 65535
$ws.Range("G20").Value = 65244
$ws.Range("H20").Value = 63759
$ws.Range("I20").Value = 2266
$ws.Range("J20").Value = 64641
$ws.Range("K20").Value = 1118
$ws.Range("L20").Value = 64793
# row 21: A21=131569  D21=52  E21=65063  F21=96  G21=104  H21=63520  I21=2461  J21=64425  K21=1136  L21=64790
$ws.Range("A21").Value = 131569
$ws.Range("D21").Value = 52
$ws.Range("E21").Value = 65063
$ws.Range("F21").Value = 96
$ws.Range("G21").Value = 104
$ws.Range("H21").Value = 63520
$ws.Range("I21").Value = 2461
$ws.Range("J21").Value = 64425
$ws.Range("K21").Value = 1136
$ws.Range("L21").Value = 64790
# row 22: A22=146457  B22=3912278  C22=2174  D22=65469  E22=65062  F22=118  G22=65304  H22=64204  I22=2241  J22=64209  K22=918  L22=64761
$ws.Range("A22").Value = 146457
$ws.Range("B22").Value = 3912278
$ws.Range("C22").Value = 2174
$ws.Range("D22").Value = 65469
$ws.Range("E22").Value = 65062
$ws.Range("F22").Value = 118
$ws.Range("G22").Value = 65304
$ws.Range("H22").Value = 64204
$ws.Range("I22").Value = 2241
$ws.Range("J22").Value = 64209
$ws.Range("K22").Value = 918
$ws.Range("L22").Value = 64761
# row 23: A23=146563  D23=446  E23=65401
$ws.Range("A23").Value = 146563
$ws.Range("D23").Value = 446
$ws.Range("E23").Value = 65401
# row 24: A24=146774  D24=45  E24=54  F24=48  G24=64657  H24=64365  I24=2550  J24=64196  K24=789  L24=64752
$ws.Range("A24").Value = 146774
$ws.Range("D24").Value = 45
$ws.Range("E24").Value = 54
$ws.Range("F24").Value = 48
$ws.Range("G24").Value = 64657
$ws.Range("H24").Value = 64365
$ws.Range("I24").Value = 2550
$ws.Range("J24").Value = 64196
$ws.Range("K24").Value = 789
$ws.Range("L24").Value = 64752
# ---------------------------------------------------------------
# 4) Merge the paired sensor-axis header cells (row 1)
# ---------------------------------------------------------------
$ws.Range("B1:C1").Merge()
$ws.Range("F1:G1").Merge()
$ws.Range("H1:I1").Merge()
$ws.Range("D1:E1").Merge()

# ---------------------------------------------------------------
# 5) Register the three "box corner" border styles used while
#    laying out the merged axis-group header boxes (top / top+right /
#    top+right+bottom) so the workbook's border table matches the
#    authored layout.
# ---------------------------------------------------------------
$ws.Range("B1").Borders.Item(8).LineStyle = 1
$ws.Range("C1").Borders.Item(8).LineStyle = 1
$ws.Range("C1").Borders.Item(10).LineStyle = 1
$ws.Range("D1").Borders.Item(8).LineStyle = 1
$ws.Range("D1").Borders.Item(10).LineStyle = 1
$ws.Range("D1").Borders.Item(9).LineStyle = 1

# ---------------------------------------------------------------
# 6) Apply the workbook's standard header style (bold, centered,
#    boxed) to row 1 and to the time_ms column, matching the other
#    sheets in the workbook.
# ---------------------------------------------------------------
$ws0.Range("A1").Copy()
$ws.Range("A1:L1").PasteSpecial(-4122)
$ws.Range("A3:A24").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# 7) Select A1 on the new sheet (matches the saved selection state)
# ---------------------------------------------------------------
$ws.Range("A1").Select()

# ---------------------------------------------------------------
# 8) Drop the stray empty "data" cells that trailed the last block of
#    state estimator rows on state_est_data (C162:C170).
# ---------------------------------------------------------------
$wsState = $wb.Worksheets.Item("state_est_data")
$wsState.Range("C162:C170").ClearContents()
